$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "solution" link for row 8 (Reverse a Linked List) in column E,
# matching the pattern used by the other rows (C/E columns hold hyperlinked
# text whose display value is the URL itself).
$ws.Range("E8").Value = "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoLibrary/LinkedList/ReverseALinkedList.cs"
$ws.Hyperlinks.Add($ws.Range("E8"), "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoLibrary/LinkedList/ReverseALinkedList.cs")
$ws.Range("E8").Style = "Hyperlink"

# Move the active selection to the newly-filled cell (also brings it into
# view, dropping the old scrolled-away topLeftCell).
$ws.Range("E8").Select() | Out-Null
